$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 41, shifting rows 41:111 down to 42:112
$ws.Rows.Item(41).Insert()

# Populate the new row 41 with the new weekly data point
$ws.Cells.Item(41, 1).Value = 4
$ws.Cells.Item(41, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(41, 3).Value = "Los Lagos"
$ws.Cells.Item(41, 4).Value = 44883
$ws.Cells.Item(41, 5).Value = 10
$ws.Cells.Item(41, 6).Value = 100112031
$ws.Cells.Item(41, 7).Value = "Poroto verde"
$ws.Cells.Item(41, 8).Value = "Magnum"
$ws.Cells.Item(41, 9).Value = "Primera"
$ws.Cells.Item(41, 10).Value = 35
$ws.Cells.Item(41, 11).Value = 36000
$ws.Cells.Item(41, 12).Value = 36000
$ws.Cells.Item(41, 13).Value = 36000
$ws.Cells.Item(41, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(41, 15).Value = "Perú"
$ws.Cells.Item(41, 16).Value = 1440
$ws.Cells.Item(41, 17).Value = 25
$ws.Cells.Item(41, 18).Value = "Hortaliza"

# Copy the date cell style/format from the row above so the new date cell
# keeps the same numeric date formatting as the rest of column D
$ws.Range("D40").Copy()
$ws.Range("D41").PasteSpecial(-4122) # xlPasteFormats
